$p = $ppt.ActivePresentation

# --- Slide 2 ("About"): split the long run to insert "(ALM) " before "solution" ---
$s2 = $p.Slides.Item(2)
$shp2 = $s2.Shapes.Item(2)
$tr2 = $shp2.TextFrame.TextRange
$full2 = $tr2.Text
$idx2 = $full2.IndexOf("solution")
$sub2 = $tr2.Characters($idx2 + 1, "solution".Length)
[void]$sub2.InsertBefore("(ALM) ")

# --- Slide 3 ("Main features (1)"): merge "   Document " + "sharing (...)" runs ---
$s3 = $p.Slides.Item(3)
$shp3 = $s3.Shapes.Item(3)
$tr3 = $shp3.TextFrame.TextRange
$full3 = $tr3.Text
$merged = "   Document sharing (including versioning and locking capabilities)"
$idx3 = $full3.IndexOf($merged)
$sub3 = $tr3.Characters($idx3 + 1, $merged.Length)
$sub3.Text = $merged
